# Update the "Fitness" (column C) values for rows 2-185 on the active sheet.
# These correspond to run_23 log data being replaced with a newer run's
# fitness progression while leaving Run/Generation columns (A/B) and rows
# 186 onward untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C4").Value = 11033
$ws.Range("C5:C9").Value = 10424
$ws.Range("C10:C15").Value = 9842
$ws.Range("C16:C19").Value = 9163
$ws.Range("C20:C25").Value = 8944
$ws.Range("C26:C52").Value = 8819
$ws.Range("C53:C63").Value = 8671
$ws.Range("C64:C80").Value = 8462
$ws.Range("C81:C106").Value = 8223
$ws.Range("C107").Value = 8049
$ws.Range("C108:C136").Value = 7750
$ws.Range("C137:C185").Value = 7295
